# Mark the remaining five vocabulary rows (18-22) as "Processed" in column C,
# matching the status already recorded for the rows above them.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C18:C22").Value = "Processed"
